$d = $word.ActiveDocument

$targetXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F"><w:r><w:t>{</w:t></w:r><w:r w:rsidR="00DE6D5A"><w:t>m</w:t></w:r><w:r w:rsidR="002033E1"><w:t>:</w:t></w:r><w:r w:rsidR="002033E1" w:rsidRPr="002033E1"><w:t xml:space="preserve"> ('</w:t></w:r><w:r w:rsidR="00D62429" w:rsidRPr="00D62429"><w:t>&lt;img src=&quot;../images/logo_M2Doc.png&quot; alt=&quot;&quot; height=&quot;54&quot;&gt;</w:t></w:r><w:r w:rsidR="002033E1" w:rsidRPr="002033E1"><w:t>&lt;h2 id=&quot;starting-with</w:t></w:r><w:r w:rsidR="00D62429"><w:t>-m2doc&quot;&gt;Starting with ' + self.</w:t></w:r><w:r w:rsidR="002033E1" w:rsidRPr="002033E1"><w:t>n</w:t></w:r><w:r w:rsidR="00D62429"><w:t>a</w:t></w:r><w:r w:rsidR="002033E1" w:rsidRPr="002033E1"><w:t>me + '&lt;/h2&gt;').from</w:t></w:r><w:r w:rsidR="00342B27"><w:t>HTML</w:t></w:r><w:r w:rsidR="00256E67"><w:t>Body</w:t></w:r><w:r w:rsidR="002033E1" w:rsidRPr="002033E1"><w:t>String(</w:t></w:r><w:r w:rsidR="00D62429"><w:t>'</w:t></w:r><w:r w:rsidR="00D62429" w:rsidRPr="00D62429"><w:t>http</w:t></w:r><w:r w:rsidR="002C5335"><w:t>s</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidR="00D62429" w:rsidRPr="00D62429"><w:t>://www.m2doc.org/tests/</w:t></w:r><w:r w:rsidR="00D62429"><w:t>'</w:t></w:r><w:r w:rsidR="002033E1" w:rsidRPr="002033E1"><w:t>)</w:t></w:r><w:r><w:t xml:space="preserve">}</w:t></w:r></w:p>
'@

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        $p.Range.InsertXML($targetXml)
        break
    }
}
